$d = $word.ActiveDocument

# Step 1: Replace all 7 occurrences of "نور چشم" with "نور والدین"
# (the chatbot was renamed from "Noor-e-Cheshm" to "Noor-e-Walidain").
# None of these matches include a quote character, so this replacement
# is unaffected by any smart-quote substitution performed during
# Find/Replace.
$result1 = $d.Content.Find.Execute('نور چشم', $true, $true, $false, $false, $false, $true, 1, $false, 'نور والدین', 2)
Write-Output "Step 1 (global rename) executed: $result1"

# Step 2: One specific instance also needs a space inserted right before
# the closing quote mark that follows the renamed text:
#   ...نور والدین" می‌توانید...   ->   ...نور والدین "‌ می‌توانید...
# Using Find.Execute's ReplaceWith parameter with a literal '"' triggers
# this runtime's automatic smart-quote substitution (a straight " turns
# into a curly “ / ” character), which would corrupt the document.
# Instead we locate the target Range with Find (empty ReplaceWith, so no
# substitution happens) and then assign its .Text property directly,
# which performs a literal, unformatted replacement that keeps the
# straight quote intact.
$rng = $d.Content
$found2 = $rng.Find.Execute('نور والدین" می‌توانید این ویدیو را مجدداً ببینید.', $true, $true, $false, $false, $false, $true, 1, $false, '', 0)
Write-Output "Step 2 target found: $found2"
if ($found2) {
    $rng.Text = 'نور والدین " می‌توانید این ویدیو را مجدداً ببینید.'
    Write-Output "Step 2 (space before quote) applied"
}
